$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.414.47"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.701.08"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.42"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5506"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2744"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06477"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.10"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07694"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.691.75"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.555"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5848"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008387"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.500.83"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.954"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.60"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.267"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.93"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1326"
$ws.Range("E25").Value = "  +7.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.936"
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06308"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.611"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.687"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.045"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6193"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.741"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01650"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.121.05"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.180"
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8854"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.31"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.80"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000109"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.260"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05277"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.137"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4306"
$ws.Range("E51").Value = "  -0.02%  "
